$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C2 was a text value "White" -> becomes numeric 8523635
$ws.Range("C2").Value = 8523635

# C3 was a text value "PollWhite" -> becomes numeric 80080
$ws.Range("C3").Value = 80080

# Move the active selection from C3 to C2
$ws.Range("C2").Select()
